# Append 5 new time-tracking entries (rows 42-46) to Sheet1, matching the
# existing table layout: Project (shared string), ProjeDate (date), Duration
# (time-of-day fraction).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 42; Project = "SpotifyPlayGen";  Date = 43334; Duration = 0.025208333333333333 },
    @{ Row = 43; Project = "ZeroToDeepLearn"; Date = 43338; Duration = 0.02783564814814815 },
    @{ Row = 44; Project = "FindYourJoy";     Date = 43340; Duration = 0.01892361111111111 },
    @{ Row = 45; Project = "CUDAproject";     Date = 43344; Duration = 0.026041666666666668 },
    @{ Row = 46; Project = "SpotifyPlayGen";  Date = 43346; Duration = 0.03349537037037037 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Project

    $ws.Cells.Item($rowNum, 2).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).NumberFormat = "m/d/yy"

    $ws.Cells.Item($rowNum, 3).Value = $r.Duration
    $ws.Cells.Item($rowNum, 3).NumberFormat = "h:mm:ss"
}

# Update the view's selection/scroll to match post-edit state (last row
# selected, similar to how Excel leaves the cursor after data entry).
[void]$ws.Range("A22").Select()
[void]$ws.Range("A47").Select()
